# Apply "Add data for 2022-11-21" update to carjacking-by-month-yoy-latest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-13"

# Update the row label for November to reflect the new "through" date
$ws.Range("A12").Value = "November (through 11-13)"

# Update the November row (row 12) values for each year column (B-I)
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 32
$ws.Range("D12").Value = 59
$ws.Range("E12").Value = 28
$ws.Range("F12").Value = 22
$ws.Range("G12").Value = 81
$ws.Range("H12").Value = 89
$ws.Range("I12").Value = 42

# Update the Total row (row 13) values for each year column (B-I)
$ws.Range("B13").Value = 273
$ws.Range("C13").Value = 518
$ws.Range("D13").Value = 769
$ws.Range("E13").Value = 643
$ws.Range("F13").Value = 504
$ws.Range("G13").Value = 1138
$ws.Range("H13").Value = 1530
$ws.Range("I13").Value = 1440
